# The post previously in row 564 ("「この子カンガルーは抱っこが好き」") was
# removed from the source data. Delete that entire row so every row below it
# (565-751) shifts up by one, and the sheet's used range shrinks from
# A1:C751 to A1:C750.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(564).Delete()
